# Rename the worksheet "Property1" -> "DataNode"
# (commit: "unify the conception of DataNode, DataTable, Entity.")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"
